$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B1").Value = 47
$ws.Range("B2").Value = 70
$ws.Range("B3").Value = 85
$ws.Range("B4").Value = 101
$ws.Range("B5").Value = 124
$ws.Range("B6").Value = 153
$ws.Range("B7").Value = 224
